# Update the Sestola report worksheet with daily COVID figures through 2022-01-05
# (rows 465-491), matching the "aggiornamento fino a 6 gennaio 2022" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, date serial (col A), nuovi pos. (col B),
# somma mobile 7gg. (col C), somma mobile 7gg. per 100mila abitanti (col D)
$data = @(
    @(465, 44539, 4, 12, 489.5960832313341),
    @(466, 44540, 0, 9, 367.1970624235006),
    @(467, 44541, 0, 8, 326.3973888208894),
    @(468, 44542, 3, 11, 448.796409628723),
    @(469, 44543, 0, 9, 367.1970624235006),
    @(470, 44544, 0, 8, 326.3973888208894),
    @(471, 44545, 0, 7, 285.5977152182783),
    @(472, 44546, 1, 4, 163.1986944104447),
    @(473, 44547, 0, 4, 163.1986944104447),
    @(474, 44548, 0, 4, 163.1986944104447),
    @(475, 44550, 0, 1, 40.79967360261118),
    @(476, 44551, 1, 2, 81.59934720522236),
    @(477, 44552, 0, 2, 81.59934720522236),
    @(478, 44553, 0, 2, 81.59934720522236),
    @(479, 44554, 3, 4, 163.1986944104447),
    @(480, 44555, 1, 5, 203.9983680130559),
    @(481, 44556, 3, 8, 326.3973888208894),
    @(482, 44557, 2, 10, 407.9967360261118),
    @(483, 44558, 6, 15, 611.9951040391677),
    @(484, 44559, 3, 18, 734.3941248470012),
    @(485, 44560, 0, 18, 734.3941248470012),
    @(486, 44561, 2, 17, 693.59445124439),
    @(487, 44562, 8, 24, 979.1921664626682),
    @(488, 44563, 1, 22, 897.592819257446),
    @(489, 44564, 3, 23, 938.3924928600571),
    @(490, 44565, 2, 19, 775.1937984496124),
    @(491, 44566, 7, 23, 938.3924928600571)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Carry the date-column formatting (style index used by column A, e.g. A464)
# down onto the newly added date cells so they match the existing column style.
$ws.Range("A464").Copy()
$ws.Range("A465:A491").PasteSpecial(-4122)
$excel.CutCopyMode = 0
